$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new column header "VIMMP_DEF" in F1, matching the style of the
# existing header cells (bold, centered, bordered).
$ws.Range("F1").Value = "VIMMP_DEF"
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# New column F values: '[]' for every mapped row except row 13
# (the "Concept" row), which carries the definition text.
$defs = @{
    2  = "[]"
    3  = "[]"
    4  = "[]"
    5  = "[]"
    6  = "[]"
    7  = "[]"
    8  = "[]"
    9  = "[]"
    10 = "[]"
    11 = "[]"
    12 = "[]"
    13 = "['An idea or notion; a unit of thought.']"
    14 = "[]"
    15 = "[]"
    16 = "[]"
    17 = "[]"
    18 = "[]"
    19 = "[]"
    20 = "[]"
}

foreach ($row in $defs.Keys) {
    $ws.Cells.Item($row, 6).Value = $defs[$row]
}
